$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 1.640985608100891
$ws.Range("B1").Value = 1.693534135818481
$ws.Range("C1").Value = 1.659261107444763
$ws.Range("D1").Value = 1.983187317848206
$ws.Range("E1").Value = 2.793979406356812
